# Insert a new weekly price record for "Cebollín" (Vega Monumental Concepción)
# as row 93, pushing the existing rows 93-128 down to 94-129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93 (rows below shift down by one,
# formatting is inherited from the row above, matching row 92's style).
$ws.Rows(93).Insert()

# Populate the newly inserted row 93 with the new record's data.
$ws.Range("A93").Value = 11
$ws.Range("B93").Value = "Vega Monumental Concepción"
$ws.Range("C93").Value = "Bíobío"
$ws.Range("D93").Value = 45135
$ws.Range("E93").Value = 8
$ws.Range("F93").Value = 100112037
$ws.Range("G93").Value = "Cebollín"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 50
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = 5000
$ws.Range("N93").Value = "`$/paquete 36 unidades"
$ws.Range("O93").Value = "Región Metropolitana"
$ws.Range("P93").Value = 139
$ws.Range("Q93").Value = 36
$ws.Range("R93").Value = "Hortaliza"
